$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new Q8 header column (J1), matching style of existing header cells
$ws.Range("I1").Copy()
$ws.Range("J1").PasteSpecial(-4122)
$ws.Range("J1").Value = "Q8"

# Update naive error values (re-simulated rt_data) and extend rows 4 and 8
$ws.Range("B2").Value = -1.726030308970723
$ws.Range("C2").Value = -0.6707712956971672
$ws.Range("D2").Value = -0.4524812228389575
$ws.Range("E2").Value = 0.6382515747836923
$ws.Range("F2").Value = 1.054974892958893
$ws.Range("G2").Value = -0.5298079032524352
$ws.Range("H2").Value = 0.4174100604179617
$ws.Range("B3").Value = -0.5368479019170082
$ws.Range("C3").Value = -0.3185578290587984
$ws.Range("D3").Value = 0.7721749685638514
$ws.Range("E3").Value = 1.188898286739052
$ws.Range("F3").Value = -0.3958845094722762
$ws.Range("G3").Value = 0.5513334541981209
$ws.Range("B4").Value = -0.2429041502562597
$ws.Range("C4").Value = 0.8478286473663901
$ws.Range("D4").Value = 1.264551965541591
$ws.Range("E4").Value = -0.3202308306697375
$ws.Range("F4").Value = 0.6269871330006596
$ws.Range("G4").Value = 0.5472930368032309
$ws.Range("H4").Value = 0.2587483631461737
$ws.Range("I4").Value = 0.2492070382450521
$ws.Range("J4").Value = -0.7701202119308102
$ws.Range("B5").Value = 1.529711793429503
$ws.Range("C5").Value = 1.946435111604704
$ws.Range("D5").Value = 0.3616523153933755
$ws.Range("E5").Value = 1.308870279063773
$ws.Range("F5").Value = 1.229176182866344
$ws.Range("G5").Value = 0.9406315092092867
$ws.Range("H5").Value = 0.9310901843081651
$ws.Range("I5").Value = -0.0882370658676972
$ws.Range("B6").Value = 1.401323249731339
$ws.Range("C6").Value = -0.1834595464799889
$ws.Range("D6").Value = 0.7637584171904082
$ws.Range("E6").Value = 0.6840643209929795
$ws.Range("F6").Value = 0.3955196473359223
$ws.Range("G6").Value = 0.3859783224348007
$ws.Range("H6").Value = -0.6333489277410616
$ws.Range("B7").Value = -0.1837880469139236
$ws.Range("C7").Value = 0.7634299167564733
$ws.Range("D7").Value = 0.6837358205590448
$ws.Range("E7").Value = 0.3951911469019876
$ws.Range("F7").Value = 0.385649822000866
$ws.Range("G7").Value = -0.6336774281749964
$ws.Range("B8").Value = 0.8767241794531259
$ws.Range("C8").Value = 0.7970300832556974
$ws.Range("D8").Value = 0.5084854095986401
$ws.Range("E8").Value = 0.4989440846975185
$ws.Range("F8").Value = -0.5203831654783438
$ws.Range("G8").Value = -0.4985743480898549
$ws.Range("H8").Value = 1.040395882322164
$ws.Range("I8").Value = 0.08839554212067069
$ws.Range("B9").Value = 0.5133193858515086
$ws.Range("C9").Value = 0.2247747121944514
$ws.Range("D9").Value = 0.2152333872933298
$ws.Range("E9").Value = -0.8040938628825325
$ws.Range("F9").Value = -0.7822850454940435
$ws.Range("G9").Value = 0.7566851849179757
$ws.Range("H9").Value = -0.195315155283518
$ws.Range("B10").Value = 0.0144045361508148
$ws.Range("C10").Value = 0.004863211249693217
$ws.Range("D10").Value = -1.014464038926169
$ws.Range("E10").Value = -0.9926552215376802
$ws.Range("F10").Value = 0.546315008874339
$ws.Range("G10").Value = -0.4056853313271546
$ws.Range("B11").Value = 0.007661802773924347
$ws.Range("C11").Value = -1.011665447401938
$ws.Range("D11").Value = -0.989856630013449
$ws.Range("E11").Value = 0.5491136003985702
$ws.Range("F11").Value = -0.4028867398029234
$ws.Range("B12").Value = -1.10825169089279
$ws.Range("C12").Value = -1.086442873504301
$ws.Range("D12").Value = 0.4525273569077186
$ws.Range("E12").Value = -0.499472983293775
$ws.Range("B13").Value = -1.011923540525774
$ws.Range("C13").Value = 0.5270466898862451
$ws.Range("D13").Value = -0.4249536503152485
$ws.Range("B14").Value = 0.58841169826065
$ws.Range("C14").Value = -0.3635886419408436
$ws.Range("B15").Value = -0.3509547008554236
